$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 1: human-readable column headers
$ws.Range("A1").Value = "Territorio:"
$ws.Range("B1").Value = "CASE  WHEN Número Trabajadores (empr) < 10 THEN '01 De 1 a 9 afiliados' WHEN Número Trabajadores (empr) < 20 THEN '02 De 10 a 19 afiliados' WHEN Número Trabajadores (empr) < 50 THEN '03 De 20 a 49 afiliados' WHEN Empre"
$ws.Range("C1").Value = "Estrato"
$ws.Range("D1").Value = "Número Empresas"
$ws.Range("E1").Value = "Dirección provincial nombre"
$ws.Range("F1").Value = "Mes y año"
$ws.Range("G1").Value = "Dirección provincial (código)"

# Row 2: measure identifiers, reordered to match row 1
$ws.Range("A2").Value = "iaest-measure:territorio"
$ws.Range("B2").Value = "iaest-measure:case--when-numero-trabajadores-empr--10-then-01-de-1-a-9-afiliados-when-numero-trabajadores-empr--20-then-02-de-10-a-19-afiliados-when-numero-trabajadores-empr--50-then-03-de-20-a-49-afiliados-when-empre"
$ws.Range("C2").Value = "iaest-measure:estrato"
$ws.Range("D2").Value = "iaest-measure:numero-empresas"
$ws.Range("E2").Value = "iaest-measure:direccion-provincial-nombre"
$ws.Range("F2").Value = "iaest-measure:mes-y-ano"
$ws.Range("G2").Value = "null"

# Row 3: "medida" for most columns, except the territorial-scope column (G) which is "null"
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "medida"
$ws.Range("D3").Value = "medida"
$ws.Range("E3").Value = "medida"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "null"

# Row 4: datatypes - numero-empresas (D) is xsd:int, rest are xsd:string, G (null) stays null
$ws.Range("A4").Value = "xsd:string"
$ws.Range("B4").Value = "xsd:string"
$ws.Range("C4").Value = "xsd:string"
$ws.Range("D4").Value = "xsd:int"
$ws.Range("E4").Value = "xsd:string"
$ws.Range("F4").Value = "xsd:string"
$ws.Range("G4").Value = "null"
